$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (so numeric-looking strings like "7.19"
# or "525.15" are not silently auto-converted to numbers by Excel's
# type inference), then strip the temporary Text number-format so the
# cell's style index is left exactly as it was before (no stray style
# bleed into the saved workbook).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "58.325.36"
$ws.Range("E2").Value = "  -2.44%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.135.75"
$ws.Range("E3").Value = "  -4.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "523.72"
$ws.Range("E5").Value = "  -4.89%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "134.35"
$ws.Range("E6").Value = "  -4.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.136.63"
$ws.Range("E8").Value = "  -4.08%  "

# Row 9 - XRP
Set-TextValue $ws.Range("D9") "0.442"
$ws.Range("E9").Value = "  -4.60%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "7.19"
$ws.Range("E10").Value = "  -7.09%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -8.30%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.379"
$ws.Range("E12").Value = "  -6.22%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "3.676.31"
$ws.Range("E13").Value = "  -4.15%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.91%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "25.47"
$ws.Range("E15").Value = "  -3.92%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "3.140.79"
$ws.Range("E16").Value = "  -4.14%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "58.315.98"
$ws.Range("E17").Value = "  -2.65%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.0000152"
$ws.Range("E18").Value = "  -6.17%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "5.75"
$ws.Range("E19").Value = "  -5.17%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "12.96"
$ws.Range("E20").Value = "  -5.75%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "7.87"
$ws.Range("E21").Value = "  -7.06%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "343.52"
$ws.Range("E22").Value = "  -7.48%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.04%  "

# Row 24 / Row 25 - swapped: Polygon/Litecoin traded places
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D24") "68.08"
$ws.Range("E24").Value = "  -7.40%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D25") "0.507"
$ws.Range("E25").Value = "  -4.28%  "

# Row 26 - WrappedeETH
Set-TextValue $ws.Range("D26") "3.276.79"
$ws.Range("E26").Value = "  -4.11%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +0.53%  "

# Row 28 - PEPE
Set-TextValue $ws.Range("D28") "0.0₃0948"
$ws.Range("E28").Value = "  -5.49%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.34%  "

# Row 30 - RenderToken
Set-TextValue $ws.Range("D30") "6.79"
$ws.Range("E30").Value = "  -3.20%  "

# Row 31 - USDe
$ws.Range("E31").Value = "  -0.06%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -7.64%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "6.87"
$ws.Range("E33").Value = "  -7.32%  "

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  -1.16%  "

# Row 35 - EthereumClassic
Set-TextValue $ws.Range("D35") "21.31"
$ws.Range("E35").Value = "  -4.78%  "

# Row 36 - NEARProtocol
Set-TextValue $ws.Range("D36") "4.78"
$ws.Range("E36").Value = "  -5.01%  "

# Row 37 - Monero
Set-TextValue $ws.Range("D37") "157.42"
$ws.Range("E37").Value = "  -4.89%  "

# Row 38 - Aptos
Set-TextValue $ws.Range("D38") "6.21"
$ws.Range("E38").Value = "  -5.69%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -9.23%  "

# Row 40 / Row 41 - swapped: Hedera/RenzoRestakedETH traded places
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D40") "3.168.50"
$ws.Range("E40").Value = "  -4.06%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D41") "0.0683"
$ws.Range("E41").Value = "  -5.45%  "

# Row 42 - OKB
Set-TextValue $ws.Range("D42") "40.42"
$ws.Range("E42").Value = "  -2.89%  "

# Row 43 - EnergySwap
Set-TextValue $ws.Range("D43") "24.09"
$ws.Range("E43").Value = "  -7.53%  "

# Row 44 / Row 45 - swapped: ONDO/Mantle traded places
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D44") "0.691"
$ws.Range("E44").Value = "  -6.88%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D45") "1.08"
$ws.Range("E45").Value = "  -1.37%  "

# Row 46 - Filecoin
Set-TextValue $ws.Range("D46") "3.89"
$ws.Range("E46").Value = "  -5.01%  "

# Row 47 - FirstDigitalUSD
$ws.Range("E47").Value = "  +0.02%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -7.73%  "

# Row 49 - Maker
Set-TextValue $ws.Range("D49") "2.272.31"
$ws.Range("E49").Value = "  -2.20%  "

# Row 50 - Cosmos
Set-TextValue $ws.Range("D50") "6.19"
$ws.Range("E50").Value = "  -2.20%  "

# Row 51 - InjectiveProtocol
Set-TextValue $ws.Range("D51") "20.73"
$ws.Range("E51").Value = "  -1.35%  "
